# Automatische test-sync: 2025-06-24 21:31:50
$wb = $excel.ActiveWorkbook

$logs = $wb.Worksheets.Item("Logs")
$dashboard = $wb.Worksheets.Item("Dashboard")

# --- New row 31: "Beschadigd product ontvangen" ---
$logs.Cells.Item(31, 1).Value = "Beschadigd product ontvangen"
$logs.Cells.Item(31, 2).Value = "mailmind.test@zohomail.eu"
$logs.Cells.Item(31, 3).Value = "Het product dat ik heb ontvangen is beschadigd aangekomen."
$logs.Cells.Item(31, 4).Value = "Retour / Terugbetaling"
$logs.Cells.Item(31, 5).Value = "Beste klant,`nBedankt voor uw bericht. Wat vervelend om te horen dat het product beschadigd is aangekomen. Kunt u ons meer details geven over de schade? Bijvoorbeeld, wat voor product is het en wat voor schade is er precies aan? Eventuele foto's van de schade kunnen ook nuttig zijn.`nZodra we deze informatie hebben, zullen we ons best doen om een passende oplossing voor u te vinden.`nMet vriendelijke groet,`n[Naam bedrijf] E-mailassistent"
$logs.Cells.Item(31, 6).Value = "2025-06-24 21:31:18"
$logs.Cells.Item(31, 7).Value = "Ja"
$logs.Rows.Item(31).AutoFit()

# --- New row 32: "Verzoek om factuur" ---
$logs.Cells.Item(32, 1).Value = "Verzoek om factuur"
$logs.Cells.Item(32, 2).Value = "mailmind.test@zohomail.eu"
$logs.Cells.Item(32, 3).Value = "Kunt u mij een factuur sturen voor mijn laatste bestelling?"
$logs.Cells.Item(32, 4).Value = "Factuur / Administratie"
$logs.Cells.Item(32, 5).Value = "Beste klant,`nBedankt voor uw e-mail. Om u te kunnen helpen met het versturen van een factuur voor uw laatste bestelling, hebben wij wat extra informatie nodig. Kunt u ons alstublieft de volgende gegevens verstrekken:`n1. Uw bestelnummer`n2. De datum van uw bestelling`n3. Het e-mailadres waarnaar wij de factuur kunnen sturen`nZodra wij deze gegevens van u hebben ontvangen, zullen wij zo spoedig mogelijk de factuur voor u opstellen en toesturen.`nMet vriendelijke groet,`n[Bedrijfsnaam] E-mailassistent"
$logs.Cells.Item(32, 6).Value = "2025-06-24 21:31:20"
$logs.Cells.Item(32, 7).Value = "Ja"
$logs.Rows.Item(32).AutoFit()

# --- Extend conditional formatting ranges to cover the new rows ---
$fcD = $logs.Range("D2:D30").FormatConditions
$fcD.Item(1).ModifyAppliesToRange($logs.Range("D2:D32"))

$fcG = $logs.Range("G2:G30").FormatConditions
$fcG.Item(1).ModifyAppliesToRange($logs.Range("G2:G32"))

# --- Update Dashboard counts ---
$dashboard.Cells.Item(2, 2).Value = 14
$dashboard.Cells.Item(3, 2).Value = 4
